$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.6803049333392641
$ws.Range("B1").Value = 4.1574160576240322
$ws.Range("C1").Value = 132758

$ws.Range("A4").Value = 2.6314883984528774
$ws.Range("B4").Value = 4.3197514589307975
$ws.Range("C4").Value = 221868
